# Generate Report for Handback
#
# Refresh the handback timestamps / priority for the files whose generated
# report rows previously shared identical values (80b09619-... in row 3 and
# ec72b5fd-... in row 5 both reported the same datetimes/priority before
# this run, so both rows move to the freshly generated values).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (col G)
$wsOverview.Range("G3").Value = "2016-09-05 08:21:28"
$wsOverview.Range("G5").Value = "2016-09-05 08:21:28"

# zh-cn sheet
$wsZhCn.Range("E3").Value = "mt"                          # Priority
$wsZhCn.Range("E5").Value = "mt"                          # Priority
$wsZhCn.Range("H3").Value = "2016-09-05 08:21:21"          # Correspond Handoff Datetime
$wsZhCn.Range("H5").Value = "2016-09-05 08:21:21"          # Correspond Handoff Datetime
$wsZhCn.Range("K3").Value = "2016-09-05 08:21:58"          # Correspond Handback DateTime
$wsZhCn.Range("K5").Value = "2016-09-05 08:21:58"          # Correspond Handback DateTime

# de-de sheet
$wsDeDe.Range("E3").Value = "mt"                          # Priority
$wsDeDe.Range("E5").Value = "mt"                          # Priority
$wsDeDe.Range("H3").Value = "2016-09-05 08:21:28"          # Correspond Handoff Datetime
$wsDeDe.Range("H5").Value = "2016-09-05 08:21:28"          # Correspond Handoff Datetime
$wsDeDe.Range("K3").Value = "2016-09-05 08:22:17"          # Correspond Handback DateTime
$wsDeDe.Range("K5").Value = "2016-09-05 08:22:17"          # Correspond Handback DateTime
